$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.904.35'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '3.116.33'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '578.39'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.45'
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.40'
$ws.Range('E9').Value = '  -3.17%  '
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '37.18'
$ws.Range('E13').Value = '  +1.90%  '
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '3.632.63'
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '66.872.54'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('D18').Value = '3.117.16'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '474.59'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.709'
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.91'
$ws.Range('E22').Value = '  +5.30%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '83.72'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.25'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.19'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.68'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '28.54'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').Value = '0.0₃0950'
$ws.Range('E33').Value = '  -6.87%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.83'
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.977'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '47.10'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.06'
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '50.18'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.311'
$ws.Range('E40').Value = '  -2.52%  '
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.60'
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').Value = '2.810.88'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '382.77'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.54'
$ws.Range('E46').Value = '  -9.77%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '135.35'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '24.86'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.20'
$ws.Range('E51').Value = '  -0.95%  '
